# Update "PERIOD TO EXPIRE" (col H) and "LAST UPDATE" (col I) for rows 3-37
# on the "Training Dashboard" sheet to reflect progress as of 04-Nov-2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$periodToExpire = @{
    3  = 350
    4  = 352
    5  = 446
    6  = 354
    7  = 470
    8  = 473
    9  = 471
    10 = 447
    11 = 471
    12 = 650
    13 = 472
    14 = 350
    15 = 475
    16 = 365
    17 = 345
    18 = 475
    19 = 473
    20 = 544
    21 = 121
    22 = 219
    23 = 219
    24 = -98
    25 = 126
    26 = 166
    27 = 126
    28 = 168
    29 = 167
    30 = 266
    31 = 267
    32 = 269
    33 = 265
    34 = 286
    35 = 287
    36 = 300
    37 = 377
}

foreach ($row in 3..37) {
    $ws.Cells.Item($row, 8).Value = $periodToExpire[$row]

    # Pre-format column I ("LAST UPDATE") as text so the date-like string
    # "04-Nov-2025" is stored as literal text (matching the source data,
    # which keeps these as plain strings) instead of being auto-converted
    # into a date serial number by Excel's smart entry detection.
    $dateCell = $ws.Cells.Item($row, 9)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "04-Nov-2025"
}
